$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 145.11111
$ws.Range("I2").Value = 163.25
$ws.Range("J2").Value = 130.6
$ws.Range("K2").Value = 163.25
$ws.Range("L2").Value = 130.6
$ws.Range("M2").Value = -50.25
$ws.Range("N2").Value = -356.6

$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()

$ws.Range("H29").Value = 6494.6665
$ws.Range("J29").Value = 6494.6665
$ws.Range("L29").Value = 19483.9995
$ws.Range("N29").Value = -20045.9995

$ws.Range("H62").Value = 2918.8
$ws.Range("I62").Value = 2854.2222
$ws.Range("K62").Value = 2854.2222
$ws.Range("M62").Value = -2230.2222

$ws.Range("H65").Value = 2918.8
$ws.Range("I65").Value = 2854.2222
$ws.Range("K65").Value = 14271.111
$ws.Range("M65").Value = -11151.111

$ws.Range("H98").Value = 2787.1538
$ws.Range("I98").Value = 2926.0833
$ws.Range("K98").Value = 2926.0833
$ws.Range("M98").Value = -1428.0833

$ws.Range("H106").Value = 6997.5
$ws.Range("I106").Value = 6995
$ws.Range("J106").Value = 7000
$ws.Range("K106").Value = 6995
$ws.Range("L106").Value = 7000
$ws.Range("M106").Value = -6364
$ws.Range("N106").Value = -8262

$ws.Range("H107").Value = 345.15
$ws.Range("J107").Value = 1405.5
$ws.Range("L107").Value = 1405.5
$ws.Range("N107").Value = -5245.5

$ws.Range("H109").Value = 49999
$ws.Range("J109").Value = 49999
$ws.Range("L109").Value = 49999
$ws.Range("N109").Value = -52773

$ws.Range("H122").Value = 2787.1538
$ws.Range("I122").Value = 2926.0833
$ws.Range("K122").Value = 8778.249899999999
$ws.Range("M122").Value = -6328.249899999999

$ws.Range("H131").Value = 1254193.2
$ws.Range("J131").Value = 5393.75
$ws.Range("L131").Value = 16181.25
$ws.Range("N131").Value = -26261.25

$ws.Range("H132").Value = 5311.212
$ws.Range("I132").Value = 5684.7856
$ws.Range("K132").Value = 17054.3568
$ws.Range("M132").Value = -14524.3568

$ws.Range("H137").Value = 1569663.1
$ws.Range("I137").Value = 10001018
$ws.Range("J137").Value = 8301.111000000001
$ws.Range("K137").Value = 30003054
$ws.Range("L137").Value = 24903.333
$ws.Range("M137").Value = -30000504
$ws.Range("N137").Value = -30003.333

$ws.Range("H141").Value = 4047.5557
$ws.Range("I141").Value = 3953.5
$ws.Range("K141").Value = 11860.5
$ws.Range("M141").Value = -6680.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3233.6775
$ws.Range("I32").Value = 2711.4424
$ws.Range("K32").Value = 2711.4424
$ws.Range("M32").Value = -2424.4424

$ws.Range("H44").Value = 67629
$ws.Range("J44").Value = 67629
$ws.Range("L44").Value = 67629
$ws.Range("N44").Value = -68605

$ws.Range("H55").Value = 76702.5
$ws.Range("J55").Value = 76702.5
$ws.Range("L55").Value = 76702.5
$ws.Range("N55").Value = -77332.5

$ws.Range("H61").Value = 4340.8335
$ws.Range("I61").Value = 2512.5
$ws.Range("J61").Value = 6169.1665
$ws.Range("K61").Value = 2512.5
$ws.Range("L61").Value = 6169.1665
$ws.Range("M61").Value = -2300.5
$ws.Range("N61").Value = -6593.1665

$ws.Range("H132").Value = 1310.6522
$ws.Range("I132").Value = 933.5854
$ws.Range("K132").Value = 2800.7562
$ws.Range("M132").Value = -270.7562000000003

$ws.Range("H136").Value = 4340.8335
$ws.Range("I136").Value = 2512.5
$ws.Range("J136").Value = 6169.1665
$ws.Range("K136").Value = 7537.5
$ws.Range("L136").Value = 18507.4995
$ws.Range("M136").Value = -4987.5
$ws.Range("N136").Value = -23607.4995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 55074
$ws.Range("J35").Value = 55074
$ws.Range("L35").Value = 55074
$ws.Range("N35").Value = -55694

$ws.Range("H86").Value = 2874.9092
$ws.Range("I86").Value = 3290
$ws.Range("J86").Value = 2529
$ws.Range("K86").Value = 3290
$ws.Range("L86").Value = 2529
$ws.Range("M86").Value = -2167
$ws.Range("N86").Value = -4775

$ws.Range("H89").Value = 2874.9092
$ws.Range("I89").Value = 3290
$ws.Range("J89").Value = 2529
$ws.Range("K89").Value = 16450
$ws.Range("L89").Value = 12645
$ws.Range("M89").Value = -10834
$ws.Range("N89").Value = -23877

$ws.Range("H134").Value = 4757.517
$ws.Range("I134").Value = 4836.5713
$ws.Range("K134").Value = 14509.7139
$ws.Range("M134").Value = -11974.7139

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 8000
$ws.Range("I12").Value = 8000
$ws.Range("K12").Value = 8000
$ws.Range("M12").Value = -7830

$ws.Range("H31").Value = 4825.4053
$ws.Range("I31").Value = 3259.8
$ws.Range("K31").Value = 3259.8
$ws.Range("M31").Value = -2964.8

$ws.Range("H34").Value = 4825.4053
$ws.Range("I34").Value = 3259.8
$ws.Range("K34").Value = 3259.8
$ws.Range("M34").Value = -3057.8

$ws.Range("H99").Value = 4041.4167
$ws.Range("I99").Value = 4199.8335
$ws.Range("J99").Value = 3883
$ws.Range("K99").Value = 4199.8335
$ws.Range("L99").Value = 3883
$ws.Range("M99").Value = -2701.8335
$ws.Range("N99").Value = -6879

$ws.Range("H107").Value = 411
$ws.Range("I107").Value = 326.85715
$ws.Range("K107").Value = 326.85715
$ws.Range("M107").Value = 1593.14285

$ws.Range("H126").Value = 4041.4167
$ws.Range("I126").Value = 4199.8335
$ws.Range("J126").Value = 3883
$ws.Range("K126").Value = 12599.5005
$ws.Range("L126").Value = 11649
$ws.Range("M126").Value = -10129.5005
$ws.Range("N126").Value = -16589

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 6610.8887
$ws.Range("I25").Value = 1000
$ws.Range("J25").Value = 7312.25
$ws.Range("K25").Value = 3000
$ws.Range("L25").Value = 21936.75
$ws.Range("M25").Value = -2831
$ws.Range("N25").Value = -22274.75

$ws.Range("H30").Value = 6610.8887
$ws.Range("I30").Value = 1000
$ws.Range("J30").Value = 7312.25
$ws.Range("K30").Value = 3000
$ws.Range("L30").Value = 21936.75
$ws.Range("M30").Value = -2898
$ws.Range("N30").Value = -22140.75

$ws.Range("H80").Value = 3227
$ws.Range("I80").Value = 3153
$ws.Range("J80").Value = 3251.6667
$ws.Range("K80").Value = 9459
$ws.Range("L80").Value = 9755.000100000001
$ws.Range("M80").Value = -8523
$ws.Range("N80").Value = -11627.0001

$ws.Range("H83").Value = 3227
$ws.Range("I83").Value = 3153
$ws.Range("J83").Value = 3251.6667
$ws.Range("K83").Value = 28377
$ws.Range("L83").Value = 29265.0003
$ws.Range("M83").Value = -23697
$ws.Range("N83").Value = -38625.0003

$ws.Range("H139").Value = 2663.6155
$ws.Range("J139").Value = 3856.8333
$ws.Range("L139").Value = 11570.4999
$ws.Range("N139").Value = -21850.4999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 50001996
$ws.Range("J80").Value = 2612.3333
$ws.Range("L80").Value = 2612.3333
$ws.Range("N80").Value = -4608.3333

$ws.Range("H83").Value = 50001996
$ws.Range("J83").Value = 2612.3333
$ws.Range("L83").Value = 13061.6665
$ws.Range("N83").Value = -23045.6665

$ws.Range("H102").Value = 2058.0908
$ws.Range("I102").Value = 1736.25
$ws.Range("K102").Value = 1736.25
$ws.Range("M102").Value = -114.25

$ws.Range("H132").Value = 2600
$ws.Range("I132").Value = 2200.3125
$ws.Range("K132").Value = 6600.9375
$ws.Range("M132").Value = -4070.9375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4502.907
$ws.Range("I40").Value = 4410.7744
$ws.Range("K40").Value = 4410.7744
$ws.Range("M40").Value = -4274.7744

$ws.Range("H43").Value = 70000
$ws.Range("J43").Value = 70000
$ws.Range("L43").Value = 70000
$ws.Range("N43").Value = -70386

$ws.Range("H132").Value = 5542.1763
$ws.Range("I132").Value = 1945.75
$ws.Range("K132").Value = 5837.25
$ws.Range("M132").Value = -3307.25

$ws.Range("H136").Value = 4676
$ws.Range("I136").Value = 4983.5454
$ws.Range("J136").Value = 3999.4
$ws.Range("K136").Value = 14950.6362
$ws.Range("L136").Value = 11998.2
$ws.Range("M136").Value = -12400.6362
$ws.Range("N136").Value = -17098.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 15666.333
$ws.Range("J51").Value = 25000
$ws.Range("L51").Value = 25000
$ws.Range("N51").Value = -26020

$ws.Range("H122").Value = 12501747
$ws.Range("J122").Value = 50002200
$ws.Range("L122").Value = 150006600
$ws.Range("N122").Value = -150011500

$ws.Range("H136").Value = 200010500
$ws.Range("I136").Value = 1000000000
$ws.Range("J136").Value = 13125
$ws.Range("K136").Value = 3000000000
$ws.Range("L136").Value = 39375
$ws.Range("M136").Value = -2999997450
$ws.Range("N136").Value = -44475
